# Update the "random" attempt columns (G, J, M) for several rows and
# recompute the corresponding "avg_random" column (N) as their average,
# matching the original workbook's convention (N = AVERAGE(G, J, M)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new values for columns G (attempt1_random), J (attempt2_random),
# M (attempt3_random), and N (avg_random = average of G, J, M).
# A value of $null means "leave unchanged".
$rowUpdates = @(
    @{ Row = 2; G = $null; J = 9;     M = $null; N = 5.666666666666667 },
    @{ Row = 3; G = 5;     J = 5;     M = 6;     N = 5.333333333333333 },
    @{ Row = 4; G = 5;     J = 8;     M = 4;     N = 5.666666666666667 },
    @{ Row = 5; G = 7;     J = 4;     M = 2;     N = 4.333333333333333 },
    @{ Row = 6; G = $null; J = $null; M = 5;     N = 7 },
    @{ Row = 7; G = 4;     J = 1;     M = 6;     N = 3.666666666666667 }
)

foreach ($update in $rowUpdates) {
    $row = $update.Row

    if ($null -ne $update.G) {
        $ws.Range("G$row").Value = $update.G
    }
    if ($null -ne $update.J) {
        $ws.Range("J$row").Value = $update.J
    }
    if ($null -ne $update.M) {
        $ws.Range("M$row").Value = $update.M
    }

    $ws.Range("N$row").Value = $update.N
}
